$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows before row 1314, pushing the existing 13 rows
# (old 1314-1326, the 2021-12 week) down to 1325-1337.
$ws.Range("A1314:R1324").Insert()

# Row 1314
$ws.Cells.Item(1314, 1).Value = 9
$ws.Cells.Item(1314, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1314, 3).Value = 'Metropolitana'
$ws.Cells.Item(1314, 4).Value = 44628
$ws.Cells.Item(1314, 5).Value = 13
$ws.Cells.Item(1314, 6).Value = 100112020
$ws.Cells.Item(1314, 7).Value = 'Tomate'
$ws.Cells.Item(1314, 8).Value = 'Larga vida'
$ws.Cells.Item(1314, 9).Value = 'Extra'
$ws.Cells.Item(1314, 10).Value = 106
$ws.Cells.Item(1314, 11).Value = 21000
$ws.Cells.Item(1314, 12).Value = 22000
$ws.Cells.Item(1314, 13).Value = 21500
$ws.Cells.Item(1314, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1314, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(1314, 16).Value = 1194
$ws.Cells.Item(1314, 17).Value = 18
$ws.Cells.Item(1314, 18).Value = 'Hortaliza'

# Row 1315
$ws.Cells.Item(1315, 1).Value = 9
$ws.Cells.Item(1315, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1315, 3).Value = 'Metropolitana'
$ws.Cells.Item(1315, 4).Value = 44628
$ws.Cells.Item(1315, 5).Value = 13
$ws.Cells.Item(1315, 6).Value = 100112020
$ws.Cells.Item(1315, 7).Value = 'Tomate'
$ws.Cells.Item(1315, 8).Value = 'Larga vida'
$ws.Cells.Item(1315, 9).Value = 'Extra'
$ws.Cells.Item(1315, 10).Value = 94
$ws.Cells.Item(1315, 11).Value = 21000
$ws.Cells.Item(1315, 12).Value = 22000
$ws.Cells.Item(1315, 13).Value = 21479
$ws.Cells.Item(1315, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1315, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1315, 16).Value = 1193
$ws.Cells.Item(1315, 17).Value = 18
$ws.Cells.Item(1315, 18).Value = 'Hortaliza'

# Row 1316
$ws.Cells.Item(1316, 1).Value = 9
$ws.Cells.Item(1316, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1316, 3).Value = 'Metropolitana'
$ws.Cells.Item(1316, 4).Value = 44628
$ws.Cells.Item(1316, 5).Value = 13
$ws.Cells.Item(1316, 6).Value = 100112020
$ws.Cells.Item(1316, 7).Value = 'Tomate'
$ws.Cells.Item(1316, 8).Value = 'Larga vida'
$ws.Cells.Item(1316, 9).Value = 'Primera'
$ws.Cells.Item(1316, 10).Value = 160
$ws.Cells.Item(1316, 11).Value = 19000
$ws.Cells.Item(1316, 12).Value = 20000
$ws.Cells.Item(1316, 13).Value = 19500
$ws.Cells.Item(1316, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1316, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(1316, 16).Value = 1083
$ws.Cells.Item(1316, 17).Value = 18
$ws.Cells.Item(1316, 18).Value = 'Hortaliza'

# Row 1317
$ws.Cells.Item(1317, 1).Value = 9
$ws.Cells.Item(1317, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1317, 3).Value = 'Metropolitana'
$ws.Cells.Item(1317, 4).Value = 44628
$ws.Cells.Item(1317, 5).Value = 13
$ws.Cells.Item(1317, 6).Value = 100112020
$ws.Cells.Item(1317, 7).Value = 'Tomate'
$ws.Cells.Item(1317, 8).Value = 'Larga vida'
$ws.Cells.Item(1317, 9).Value = 'Primera'
$ws.Cells.Item(1317, 10).Value = 106
$ws.Cells.Item(1317, 11).Value = 19000
$ws.Cells.Item(1317, 12).Value = 20000
$ws.Cells.Item(1317, 13).Value = 19500
$ws.Cells.Item(1317, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1317, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1317, 16).Value = 1083
$ws.Cells.Item(1317, 17).Value = 18
$ws.Cells.Item(1317, 18).Value = 'Hortaliza'

# Row 1318
$ws.Cells.Item(1318, 1).Value = 9
$ws.Cells.Item(1318, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1318, 3).Value = 'Metropolitana'
$ws.Cells.Item(1318, 4).Value = 44628
$ws.Cells.Item(1318, 5).Value = 13
$ws.Cells.Item(1318, 6).Value = 100112020
$ws.Cells.Item(1318, 7).Value = 'Tomate'
$ws.Cells.Item(1318, 8).Value = 'Larga vida'
$ws.Cells.Item(1318, 9).Value = 'Segunda'
$ws.Cells.Item(1318, 10).Value = 97
$ws.Cells.Item(1318, 11).Value = 17000
$ws.Cells.Item(1318, 12).Value = 18000
$ws.Cells.Item(1318, 13).Value = 17495
$ws.Cells.Item(1318, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1318, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(1318, 16).Value = 972
$ws.Cells.Item(1318, 17).Value = 18
$ws.Cells.Item(1318, 18).Value = 'Hortaliza'

# Row 1319
$ws.Cells.Item(1319, 1).Value = 9
$ws.Cells.Item(1319, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1319, 3).Value = 'Metropolitana'
$ws.Cells.Item(1319, 4).Value = 44628
$ws.Cells.Item(1319, 5).Value = 13
$ws.Cells.Item(1319, 6).Value = 100112020
$ws.Cells.Item(1319, 7).Value = 'Tomate'
$ws.Cells.Item(1319, 8).Value = 'Larga vida'
$ws.Cells.Item(1319, 9).Value = 'Segunda'
$ws.Cells.Item(1319, 10).Value = 79
$ws.Cells.Item(1319, 11).Value = 17000
$ws.Cells.Item(1319, 12).Value = 18000
$ws.Cells.Item(1319, 13).Value = 17494
$ws.Cells.Item(1319, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1319, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1319, 16).Value = 972
$ws.Cells.Item(1319, 17).Value = 18
$ws.Cells.Item(1319, 18).Value = 'Hortaliza'

# Row 1320
$ws.Cells.Item(1320, 1).Value = 9
$ws.Cells.Item(1320, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1320, 3).Value = 'Metropolitana'
$ws.Cells.Item(1320, 4).Value = 44628
$ws.Cells.Item(1320, 5).Value = 13
$ws.Cells.Item(1320, 6).Value = 100112020
$ws.Cells.Item(1320, 7).Value = 'Tomate'
$ws.Cells.Item(1320, 8).Value = 'Larga vida'
$ws.Cells.Item(1320, 9).Value = 'Tercera'
$ws.Cells.Item(1320, 10).Value = 79
$ws.Cells.Item(1320, 11).Value = 15000
$ws.Cells.Item(1320, 12).Value = 16000
$ws.Cells.Item(1320, 13).Value = 15494
$ws.Cells.Item(1320, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1320, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(1320, 16).Value = 861
$ws.Cells.Item(1320, 17).Value = 18
$ws.Cells.Item(1320, 18).Value = 'Hortaliza'

# Row 1321
$ws.Cells.Item(1321, 1).Value = 9
$ws.Cells.Item(1321, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1321, 3).Value = 'Metropolitana'
$ws.Cells.Item(1321, 4).Value = 44628
$ws.Cells.Item(1321, 5).Value = 13
$ws.Cells.Item(1321, 6).Value = 100112020
$ws.Cells.Item(1321, 7).Value = 'Tomate'
$ws.Cells.Item(1321, 8).Value = 'Larga vida'
$ws.Cells.Item(1321, 9).Value = 'Tercera'
$ws.Cells.Item(1321, 10).Value = 52
$ws.Cells.Item(1321, 11).Value = 15000
$ws.Cells.Item(1321, 12).Value = 16000
$ws.Cells.Item(1321, 13).Value = 15500
$ws.Cells.Item(1321, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1321, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1321, 16).Value = 861
$ws.Cells.Item(1321, 17).Value = 18
$ws.Cells.Item(1321, 18).Value = 'Hortaliza'

# Row 1322
$ws.Cells.Item(1322, 1).Value = 9
$ws.Cells.Item(1322, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1322, 3).Value = 'Metropolitana'
$ws.Cells.Item(1322, 4).Value = 44628
$ws.Cells.Item(1322, 5).Value = 13
$ws.Cells.Item(1322, 6).Value = 100112020
$ws.Cells.Item(1322, 7).Value = 'Tomate'
$ws.Cells.Item(1322, 8).Value = 'Semiduro'
$ws.Cells.Item(1322, 9).Value = 'Primera'
$ws.Cells.Item(1322, 10).Value = 160
$ws.Cells.Item(1322, 11).Value = 12000
$ws.Cells.Item(1322, 12).Value = 13000
$ws.Cells.Item(1322, 13).Value = 12500
$ws.Cells.Item(1322, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1322, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1322, 16).Value = 694
$ws.Cells.Item(1322, 17).Value = 18
$ws.Cells.Item(1322, 18).Value = 'Hortaliza'

# Row 1323
$ws.Cells.Item(1323, 1).Value = 9
$ws.Cells.Item(1323, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1323, 3).Value = 'Metropolitana'
$ws.Cells.Item(1323, 4).Value = 44628
$ws.Cells.Item(1323, 5).Value = 13
$ws.Cells.Item(1323, 6).Value = 100112020
$ws.Cells.Item(1323, 7).Value = 'Tomate'
$ws.Cells.Item(1323, 8).Value = 'Semiduro'
$ws.Cells.Item(1323, 9).Value = 'Segunda'
$ws.Cells.Item(1323, 10).Value = 79
$ws.Cells.Item(1323, 11).Value = 10000
$ws.Cells.Item(1323, 12).Value = 11000
$ws.Cells.Item(1323, 13).Value = 10494
$ws.Cells.Item(1323, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1323, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1323, 16).Value = 583
$ws.Cells.Item(1323, 17).Value = 18
$ws.Cells.Item(1323, 18).Value = 'Hortaliza'

# Row 1324
$ws.Cells.Item(1324, 1).Value = 9
$ws.Cells.Item(1324, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1324, 3).Value = 'Metropolitana'
$ws.Cells.Item(1324, 4).Value = 44628
$ws.Cells.Item(1324, 5).Value = 13
$ws.Cells.Item(1324, 6).Value = 100112020
$ws.Cells.Item(1324, 7).Value = 'Tomate'
$ws.Cells.Item(1324, 8).Value = 'Semiduro'
$ws.Cells.Item(1324, 9).Value = 'Tercera'
$ws.Cells.Item(1324, 10).Value = 43
$ws.Cells.Item(1324, 11).Value = 8000
$ws.Cells.Item(1324, 12).Value = 9000
$ws.Cells.Item(1324, 13).Value = 8512
$ws.Cells.Item(1324, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(1324, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(1324, 16).Value = 473
$ws.Cells.Item(1324, 17).Value = 18
$ws.Cells.Item(1324, 18).Value = 'Hortaliza'
